$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1798.6842
$ws.Range("I92").Value = 420.76923
$ws.Range("J92").Value = 4784.1665
$ws.Range("K92").Value = 420.76923
$ws.Range("L92").Value = 4784.1665
$ws.Range("M92").Value = 827.23077
$ws.Range("N92").Value = -7280.1665
$ws.Range("H98").Value = 1926.091
$ws.Range("I98").Value = 2018.7
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 2018.7
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = -520.7
$ws.Range("N98").Value = -3996
$ws.Range("H122").Value = 1926.091
$ws.Range("I122").Value = 2018.7
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 6056.1
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -3606.1
$ws.Range("N122").Value = -7900
$ws.Range("H134").Value = 36980
$ws.Range("J134").Value = 36980
$ws.Range("L134").Value = 36980
$ws.Range("N134").Value = -47120
$ws.Range("H137").Value = 1402.9697
$ws.Range("I137").Value = 736
$ws.Range("K137").Value = 2208
$ws.Range("M137").Value = 342
$ws.Range("H138").Value = 3690.5833
$ws.Range("I138").Value = 837.7368
$ws.Range("J138").Value = 4524.492
$ws.Range("K138").Value = 2513.2104
$ws.Range("L138").Value = 13573.476
$ws.Range("M138").Value = 2626.7896
$ws.Range("N138").Value = -23853.476

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16671015
$ws.Range("I32").Value = 20835972
$ws.Range("K32").Value = 20835972
$ws.Range("M32").Value = -20835685
$ws.Range("H61").Value = 2839
$ws.Range("I61").Value = 1841.3334
$ws.Range("J61").Value = 4335.5
$ws.Range("K61").Value = 1841.3334
$ws.Range("L61").Value = 4335.5
$ws.Range("M61").Value = -1629.3334
$ws.Range("N61").Value = -4759.5
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494
$ws.Range("H122").Value = 1546.3334
$ws.Range("I122").Value = 1479.6
$ws.Range("K122").Value = 4438.799999999999
$ws.Range("M122").Value = -1988.799999999999
$ws.Range("H136").Value = 2839
$ws.Range("I136").Value = 1841.3334
$ws.Range("J136").Value = 4335.5
$ws.Range("K136").Value = 5524.0002
$ws.Range("L136").Value = 13006.5
$ws.Range("M136").Value = -2974.0002
$ws.Range("N136").Value = -18106.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1075.1428
$ws.Range("I37").Value = 921
$ws.Range("J37").Value = 2000
$ws.Range("K37").Value = 921
$ws.Range("L37").Value = 2000
$ws.Range("M37").Value = -784
$ws.Range("N37").Value = -2274
$ws.Range("H94").Value = 781.73334
$ws.Range("I94").Value = 766.1429
$ws.Range("K94").Value = 766.1429
$ws.Range("M94").Value = -315.1429000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 27500
$ws.Range("J52").Value = 27500
$ws.Range("L52").Value = 27500
$ws.Range("N52").Value = -28088
$ws.Range("H99").Value = 1563.3334
$ws.Range("I99").Value = 1552.2
$ws.Range("K99").Value = 1552.2
$ws.Range("M99").Value = -54.20000000000005
$ws.Range("H107").Value = 710.0606
$ws.Range("I107").Value = 639.55554
$ws.Range("J107").Value = 794.6667
$ws.Range("K107").Value = 639.55554
$ws.Range("L107").Value = 794.6667
$ws.Range("M107").Value = 1280.44446
$ws.Range("N107").Value = -4634.6667
$ws.Range("H126").Value = 1563.3334
$ws.Range("I126").Value = 1552.2
$ws.Range("K126").Value = 4656.6
$ws.Range("M126").Value = -2186.6
$ws.Range("H129").Value = 25999.6
$ws.Range("I129").Value = 10000
$ws.Range("J129").Value = 49999
$ws.Range("K129").Value = 10000
$ws.Range("L129").Value = 49999
$ws.Range("M129").Value = -5000
$ws.Range("N129").Value = -59999
$ws.Range("H134").Value = 1869.3864
$ws.Range("I134").Value = 1250.5312
$ws.Range("J134").Value = 3519.6667
$ws.Range("K134").Value = 3751.5936
$ws.Range("L134").Value = 10559.0001
$ws.Range("M134").Value = -1216.5936
$ws.Range("N134").Value = -15629.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 435.72726
$ws.Range("I98").Value = 428.7143
$ws.Range("J98").Value = 448
$ws.Range("K98").Value = 1286.1429
$ws.Range("L98").Value = 1344
$ws.Range("M98").Value = 211.8571000000002
$ws.Range("N98").Value = -4340
$ws.Range("H113").Value = 1379898
$ws.Range("I113").Value = 4926636.5
$ws.Range("J113").Value = 610.8333
$ws.Range("K113").Value = 14779909.5
$ws.Range("L113").Value = 1832.4999
$ws.Range("M113").Value = -14777739.5
$ws.Range("N113").Value = -6172.4999
$ws.Range("H122").Value = 6739.6855
$ws.Range("I122").Value = 10675.143
$ws.Range("K122").Value = 96076.287
$ws.Range("M122").Value = -93626.287
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = ""
$ws.Range("N123").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2931.077
$ws.Range("I80").Value = 2808.3333
$ws.Range("J80").Value = 3127.4666
$ws.Range("K80").Value = 2808.3333
$ws.Range("L80").Value = 3127.4666
$ws.Range("M80").Value = -1810.3333
$ws.Range("N80").Value = -5123.4666
$ws.Range("H83").Value = 2931.077
$ws.Range("I83").Value = 2808.3333
$ws.Range("J83").Value = 3127.4666
$ws.Range("K83").Value = 14041.6665
$ws.Range("L83").Value = 15637.333
$ws.Range("M83").Value = -9049.6665
$ws.Range("N83").Value = -25621.333
$ws.Range("H135").Value = 46980
$ws.Range("J135").Value = 46980
$ws.Range("L135").Value = 46980
$ws.Range("N135").Value = -57120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 17028.334
$ws.Range("J63").Value = 17028.334
$ws.Range("L63").Value = 17028.334
$ws.Range("N63").Value = -18526.334
$ws.Range("H66").Value = 17028.334
$ws.Range("J66").Value = 17028.334
$ws.Range("L66").Value = 51085.00199999999
$ws.Range("N66").Value = -58573.00199999999
$ws.Range("H74").Value = 66000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 66000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 66000
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = -67996
$ws.Range("H77").Value = 66000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 66000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 198000
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = -207984
$ws.Range("H93").Value = 12932.875
$ws.Range("I93").Value = 14709
$ws.Range("K93").Value = 14709
$ws.Range("M93").Value = -13461
$ws.Range("H119").Value = 45684
$ws.Range("J119").Value = 45684
$ws.Range("L119").Value = 45684
$ws.Range("N119").Value = -55360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 978.1429
$ws.Range("I81").Value = 927.5
$ws.Range("J81").Value = 998.4
$ws.Range("K81").Value = 1855
$ws.Range("L81").Value = 1996.8
$ws.Range("M81").Value = -794
$ws.Range("N81").Value = -4118.8
$ws.Range("H84").Value = 978.1429
$ws.Range("I84").Value = 927.5
$ws.Range("J84").Value = 998.4
$ws.Range("K84").Value = 9275
$ws.Range("L84").Value = 9984
$ws.Range("M84").Value = -3971
$ws.Range("N84").Value = -20592
$ws.Range("H119").Value = 42349
$ws.Range("J119").Value = 42349
$ws.Range("L119").Value = 42349
$ws.Range("N119").Value = -52025
